# "Generate Report for Handoff"
#
# The localization status report moves the zh-cn / de-de translation jobs
# from "In Translation" to "Ready for handoff", and refreshes the
# handoff-related timestamps to the moment this report was generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# Status columns for both locales flip to "Ready for handoff", and the
# "Latest HO Xliff Generate Date" timestamp is refreshed.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-25 10:40:45"

# --- zh-cn detail sheet ----------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-25 10:40:40"

# --- de-de detail sheet ----------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-25 10:40:45"

# --- Resize the Status columns to fit the new, longer text -----------
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).AutoFit()
